# newswitch.xlsx — finished get content and download page source
#
# The scraper appended the freshly-downloaded "IoT" article URLs to the
# IoT worksheet (previously just the two header rows) and left the
# workbook with the "AI" tab active/selected at A3 (the first tab,
# scrolled to the top) instead of the last tab ("VR広告").

$wb = $excel.ActiveWorkbook

# --- 1. Populate the IoT sheet with the newly scraped article URLs ----
$ws = $wb.Worksheets.Item("IoT")

$urls = @(
    "https://newswitch.jp/outline/20289",
    "https://newswitch.jp/outline/20283",
    "https://newswitch.jp/outline/20274",
    "https://newswitch.jp/outline/20260",
    "https://newswitch.jp/outline/20228",
    "https://newswitch.jp/outline/20192",
    "https://newswitch.jp/outline/20225",
    "https://newswitch.jp/outline/20173",
    "https://newswitch.jp/outline/20177",
    "https://newswitch.jp/outline/20170",
    "https://newswitch.jp/outline/20153",
    "https://newswitch.jp/outline/20148",
    "https://newswitch.jp/outline/20135",
    "https://newswitch.jp/outline/20109",
    "https://newswitch.jp/outline/19930",
    "https://newswitch.jp/outline/20058",
    "https://newswitch.jp/outline/20066",
    "https://newswitch.jp/outline/20036",
    "https://newswitch.jp/outline/20032",
    "https://newswitch.jp/outline/20020"
)

$startRow = 3
for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $urls[$i]
}

# New cells should carry the plain/default style (no fill), not the
# column's default style used by the "url_articles" label in A2.
$lastRow = $startRow + $urls.Length - 1
$ws.Range("A" + $startRow + ":A" + $lastRow).Style = "Normal"

# --- 2. Move the active tab/selection back to the first sheet ("AI") --
$ai = $wb.Worksheets.Item("AI")
[void]$ai.Activate()
[void]$ai.Range("A3").Select()
